$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6 (shifts existing rows 6-14 down to 7-15)
$ws.Rows.Item(6).Insert()

# Row 2
$ws.Range("A2").Value = 'Win32\Release\Format.exe.RCData.TFORM1.Caption'
$ws.Range("B2").Value = 'Sample'
$ws.Range("C2").Value = 'Esimerkki'
$ws.Range("E2").Value = 'Probe'
$ws.Range("F2").Value = 'Steekproef'

# Row 3
$ws.Range("A3").Value = 'Win32\Release\Format.exe.RCData.TFORM1.FirstNameLabel.Caption'
$ws.Range("B3").Value = '&First name:'
$ws.Range("C3").Value = '&Ensimmäinen nimi:'
$ws.Range("E3").Value = '&Vorname, Vorname:'
$ws.Range("F3").Value = '&Voornaam:'

# Row 4
$ws.Range("A4").Value = 'Win32\Release\Format.exe.RCData.TFORM1.CountLabel.Caption'
$ws.Range("B4").Value = '&Count:'
$ws.Range("C4").Value = '&Määrä:'
$ws.Range("E4").Value = '&Anzahl:'
$ws.Range("F4").Value = '&Telling:'

# Row 5
$ws.Range("A5").Value = 'Win32\Release\Format.exe.RCData.TFORM1.SecondNameLabel.Caption'
$ws.Range("B5").Value = '&Second name:'
$ws.Range("C5").Value = '&Toinen nimi:'
$ws.Range("E5").Value = '&Zweiter Name:'
$ws.Range("F5").Value = '&Secondenaam:'

# Row 6
$ws.Range("A6").Value = 'Win32\Release\Format.exe.RCData.TFORM1.CountEdit.Text'
$ws.Range("B6").Value = '0'

# Row 7
$ws.Range("A7").Value = 'Win32\Release\Format.exe.RCData.TFORM1.LanguageButton.Caption'
$ws.Range("B7").Value = '&Language...'
$ws.Range("C7").Value = '&Kieli...'
$ws.Range("E7").Value = '&Sprache....'
$ws.Range("F7").Value = '&Taal...'

# Row 8
$ws.Range("A8").Value = 'Win32\Release\Format.exe.RCData.TNTLANGUAGEDIALOG.Caption'
$ws.Range("B8").Value = 'Select Language'
$ws.Range("C8").Value = 'Valitse kieli'
$ws.Range("E8").Value = 'Sprache auswählen'
$ws.Range("F8").Value = 'Selecteer Taal'

# Row 9
$ws.Range("A9").Value = 'Win32\Release\Format.exe.RCData.TNTLANGUAGEDIALOG.OkButton.Caption'
$ws.Range("B9").Value = 'OK'
$ws.Range("C9").Value = 'OK'
$ws.Range("E9").Value = 'OK'
$ws.Range("F9").Value = 'OK'

# Row 10
$ws.Range("A10").Value = 'Win32\Release\Format.exe.RCData.TNTLANGUAGEDIALOG.CancelButton.Caption'
$ws.Range("B10").Value = 'Cancel'
$ws.Range("C10").Value = 'Peru'
$ws.Range("E10").Value = 'Abbrechen'
$ws.Range("F10").Value = 'Annuleren'

# Row 11
$ws.Range("A11").Value = 'Win32\Release\Format.exe.String.Unit1.SHello'
$ws.Range("B11").Value = 'Hello {0}!'
$ws.Range("C11").Value = 'Moi {0}!'
$ws.Range("E11").Value = 'Hallo {0}!'
$ws.Range("F11").Value = 'Hallo {0}!'

# Row 12
$ws.Range("A12").Value = 'Win32\Release\Format.exe.String.Unit1.SHello2'
$ws.Range("B12").Value = 'Hello {0} and {1}!'
$ws.Range("C12").Value = 'Moi {0} ja {1}!'
$ws.Range("D12").Value = 'Moi {1} ja {0}!'
$ws.Range("E12").Value = 'Hallo {0} und {1}!'
$ws.Range("F12").Value = 'Hallo {0} en {1}!'
$ws.Range("G12").Value = 'שלום {0} {1}!'

# Row 13
$ws.Range("A13").Value = 'Win32\Release\Format.exe.String.Unit1.SCount'
$ws.Range("B13").Value = '{0} has {1} cars'
$ws.Range("C13").Value = '{0}:lla on {1} autoa'
$ws.Range("E13").Value = '{0} hat {1} Autos'
$ws.Range("F13").Value = '{0} heeft {1} auto''s'

# Row 14
$ws.Range("A14").Value = 'Win32\Release\Format.exe.String.Unit1.SCount2'
$ws.Range("B14").Value = '{0} cars will pick up {1} and {2}'
$ws.Range("C14").Value = '{0} autoa hakee {1}:in ja {2}:in'
$ws.Range("E14").Value = 'Autos nehmen {1} und {2} auf.'
$ws.Range("F14").Value = 'Auto''s zullen {0} ophalen {1} en {2}.'
$ws.Range("G14").Value = '{0} מכוניות יאספו {1} ו-{2}'

# Row 15
$ws.Range("A15").Value = 'Win32\Release\Format.exe.String.Unit1.SDouble'
$ws.Range("B15").Value = '{0} swims and {0} skis'
$ws.Range("C15").Value = '{0} ui ja {0} hiihtää'
$ws.Range("D15").Value = '{0} ui ja {0} hiihtää'
$ws.Range("E15").Value = 'Schwimmen und Skier'
$ws.Range("F15").Value = '{0} Zwemmen en {0} ski''s'

